# "removing joke that did not age well"
# The deck's final slide ("Interlude" - the Mazeltov Cocktail joke slide)
# is removed entirely from the presentation.

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$slide = $p.Slides.Item($lastIndex)

$slide.Delete()
